$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.918431162834167
$ws.Range("B1").Value = 3.57275128364563
$ws.Range("C1").Value = 3.298248291015625
$ws.Range("D1").Value = 3.576974391937256
$ws.Range("E1").Value = 1.2140212059021
